$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New week header columns
$ws.Range("AO1").Value = "Week 40"
$ws.Range("AP1").Value = "Week 41"

# New data entries for week 40 / week 41
$ws.Range("AO2").Value = 10      # Scott Foxley-Berry
$ws.Range("AO4").Value = 4       # Laura Thompson
$ws.Range("AP4").Value = 5       # Laura Thompson
$ws.Range("AP5").Value = 5.5     # Jason Bohrer
$ws.Range("AO7").Value = 5.5     # Jason Liess
$ws.Range("AP8").Value = 2.75    # Kim Quan
$ws.Range("AO9").Value = 1.5     # Leo Hayward
$ws.Range("AP9").Value = 1.5     # Leo Hayward
$ws.Range("AO10").Value = 6      # Daniel Burcham
$ws.Range("AP10").Value = 4.5    # Daniel Burcham

# Update the active selection to AP9 as in the diff
$ws.Range("AP9").Select()
